$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.1
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 2.25
$ws.Range("K2").Value = 2.05
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67
$ws.Range("AE2").Value = 8.5
$ws.Range("AI2").Value = 351
$ws.Range("AJ2").Value = 7
